$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, reporting week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/17/2023  Through  4/23/2023"

# --- Weekly crime-complaint table updates (rows 14-30, 37) ---
# Row 14
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0"
$ws.Range("C14").NumberFormat = "General"
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("D14").Value = 1
$ws.Range("E14").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E14").Value = -100
$ws.Range("G14").NumberFormat = "#,##0"
$ws.Range("G14").Value = 1
$ws.Range("H14").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 6
$ws.Range("K14").Value = -83.333333333333
# Row 15
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C15").NumberFormat = "General"
$ws.Range("N15").Value = -50
# Row 16
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -22.222222222222
$ws.Range("F16").Value = 30
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 117
$ws.Range("J16").Value = 115
$ws.Range("K16").Value = 1.739130434782
$ws.Range("L16").Value = 39.285714285714
$ws.Range("M16").Value = -16.428571428571
$ws.Range("N16").Value = -75.625
# Row 17
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -11.111111111111
$ws.Range("F17").Value = 50
$ws.Range("G17").Value = 46
$ws.Range("H17").Value = 8.695652173913
$ws.Range("I17").Value = 187
$ws.Range("J17").Value = 170
$ws.Range("K17").Value = 10
$ws.Range("L17").Value = 41.666666666666
$ws.Range("M17").Value = 41.666666666666
$ws.Range("N17").Value = -0.531914893617
# Row 18
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 150
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 74
$ws.Range("J18").Value = 66
$ws.Range("K18").Value = 12.121212121212
$ws.Range("L18").Value = 27.586206896551
$ws.Range("M18").Value = -33.928571428571
$ws.Range("N18").Value = -89.290882778581
# Row 19
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 72.727272727272
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = 14.893617021276
$ws.Range("I19").Value = 193
$ws.Range("J19").Value = 203
$ws.Range("K19").Value = -4.926108374384
$ws.Range("L19").Value = 6.043956043956
$ws.Range("M19").Value = 14.201183431952
$ws.Range("N19").Value = -21.224489795918
# Row 20
$ws.Range("C20").Value = 9
$ws.Range("E20").Value = 80
$ws.Range("F20").Value = 25
$ws.Range("H20").Value = 4.166666666666
$ws.Range("I20").Value = 102
$ws.Range("J20").Value = 83
$ws.Range("K20").Value = 22.89156626506
$ws.Range("L20").Value = 96.153846153846
$ws.Range("M20").Value = 104
$ws.Range("N20").Value = -78.705636743215
# Row 21
$ws.Range("C21").Value = 48
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = 29.729729729729
$ws.Range("F21").Value = 182
$ws.Range("G21").Value = 168
$ws.Range("H21").Value = 8.333333333333
$ws.Range("I21").Value = 685
$ws.Range("J21").Value = 650
$ws.Range("K21").Value = 5.384615384615
$ws.Range("L21").Value = 32.495164410058
$ws.Range("M21").Value = 11.928104575163
$ws.Range("N21").Value = -67.673430863614
# Row 22
$ws.Range("C22").Value = 2
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 2
$ws.Range("E22").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 5
$ws.Range("G22").NumberFormat = "#,##0"
$ws.Range("G22").Value = 2
$ws.Range("H22").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("H22").Value = 150
$ws.Range("I22").Value = 15
$ws.Range("J22").Value = 12
$ws.Range("K22").Value = 25
$ws.Range("L22").Value = 87.5
$ws.Range("M22").Value = 87.5
# Row 23
$ws.Range("L23").Value = -33.333333333333
# Row 24
$ws.Range("C24").Value = 65
$ws.Range("D24").Value = 46
$ws.Range("E24").Value = 41.304347826087
$ws.Range("F24").Value = 246
$ws.Range("G24").Value = 260
$ws.Range("H24").Value = -5.384615384615
$ws.Range("I24").Value = 772
$ws.Range("J24").Value = 994
$ws.Range("K24").Value = -22.334004024144
$ws.Range("L24").Value = 104.774535809019
$ws.Range("M24").Value = 94.458438287153
# Row 25
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -13.333333333333
$ws.Range("F25").Value = 70
$ws.Range("G25").Value = 70
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 227
$ws.Range("J25").Value = 237
$ws.Range("K25").Value = -4.2194092827
$ws.Range("L25").Value = 17.61658031088
$ws.Range("M25").Value = -11.673151750972
# Row 26
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C26").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("E26").NumberFormat = "General"
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 25
# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("E27").NumberFormat = "General"
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 14.285714285714
$ws.Range("I27").Value = 29
$ws.Range("K27").Value = 16
$ws.Range("L27").Value = -6.451612903225
# Row 28
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C28").NumberFormat = "General"
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 7
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = -30
$ws.Range("M28").Value = 40
$ws.Range("N28").Value = -56.25
# Row 29
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C29").NumberFormat = "General"
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 7
$ws.Range("J29").Value = 12
$ws.Range("K29").Value = -41.666666666666
$ws.Range("L29").Value = -30
$ws.Range("M29").Value = 40
$ws.Range("N29").Value = -41.666666666666
